$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column D; this shifts existing D:K data to E:L
$ws.Columns("D").Insert()

# The newly inserted column D picks up formatting from column C (to its left).
# Copy number formats from column E (which holds the shifted original D formatting)
# into column D so the new column matches the existing data columns' styles.
$ws.Range("E5:E102").Copy() | Out-Null
$ws.Range("D5:D102").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new column D with the latest reporting period's figures
$ws.Range("D7").Value = 43372
$ws.Range("D8").Value = 1344400
$ws.Range("D9").Value = 592300
$ws.Range("D10").Value = 752100
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 125200
$ws.Range("D17").Value = 1191100
$ws.Range("D18").Value = 153300
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 278500
$ws.Range("D22").Value = 77400
$ws.Range("D23").Value = 75900
$ws.Range("D24").Value = 500
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 75400
$ws.Range("D27").Value = 75400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 1100
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 76500
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 76500
$ws.Range("D38").Value = 43372
$ws.Range("D41").Value = 5200
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 71300
$ws.Range("D44").Value = 59100
$ws.Range("D45").Value = 22200
$ws.Range("D46").Value = 157800
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 649200
$ws.Range("D49").Value = 1268700
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 25600
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2101200
$ws.Range("D57").Value = 38300
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 180800
$ws.Range("D60").Value = 219000
$ws.Range("D61").Value = 1255100
$ws.Range("D62").Value = 133200
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1607400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 493800
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43372
$ws.Range("D81").Value = 76500
$ws.Range("D83").Value = 125200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 208500
$ws.Range("D91").Value = -32900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -39100
$ws.Range("D96").Value = -147200
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -167100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 2400
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = "NA"
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"
$ws.Range("K29").Value = 0
$ws.Range("E89").Value = 161300
$ws.Range("F89").Value = 157400
$ws.Range("G89").Value = 324200
$ws.Range("H89").Value = 225600
$ws.Range("I89").Value = 214300
$ws.Range("J89").Value = 111000
$ws.Range("K89").Value = 132800
$ws.Range("E100").Value = -172900
$ws.Range("F100").Value = -218500
$ws.Range("G100").Value = -228500
$ws.Range("H100").Value = -223600
$ws.Range("I100").Value = -226700
$ws.Range("J100").Value = 113500
$ws.Range("K100").Value = -120600

